$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Weekly Quantity": append rows 25 and 26
# ---------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @{ Row = 25; A = 45662.99999999999; B = 23 },
    @{ Row = 26; A = 45669.99999999999; B = 7 }
)

foreach ($r in $weeklyNewRows) {
    $wsWeekly.Cells.Item($r.Row, 1).NumberFormat = $wsWeekly.Cells.Item($r.Row - 1, 1).NumberFormat
    $wsWeekly.Cells.Item($r.Row, 1).Value = $r.A
    $wsWeekly.Cells.Item($r.Row, 2).Value = $r.B
}

# ---------------------------------------------------------------
# Sheet "Monthly Trend": append row 13
# ---------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(13, 1).NumberFormat = $wsMonthly.Cells.Item(12, 1).NumberFormat
$wsMonthly.Cells.Item(13, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(13, 2).Value = 30

# ---------------------------------------------------------------
# Sheet "PO Forecast": new forecast model values
# ---------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Rows 2-24: only the PO_Forecast (column B) values change
$forecastUpdates = @{
    2 = 8; 3 = 9; 4 = 10; 5 = 10; 6 = 11; 7 = 12; 8 = 12; 9 = 13; 10 = 13;
    11 = 14; 12 = 14; 13 = 15; 14 = 15; 15 = 16; 16 = 16; 17 = 17; 18 = 17;
    19 = 18; 20 = 19; 21 = 21; 22 = 22; 23 = 23; 24 = 24
}

foreach ($row in $forecastUpdates.Keys) {
    $wsForecast.Cells.Item($row, 2).Value = $forecastUpdates[$row]
}

# Rows 25-34: both the date (column A) and forecast (column B) change;
# this also extends the sheet from 32 to 34 rows.
$forecastTail = @(
    @{ Row = 25; A = 45662.99999999999; B = 26 },
    @{ Row = 26; A = 45669.99999999999; B = 26 },
    @{ Row = 27; A = 45676.99999999999; B = 27 },
    @{ Row = 28; A = 45683.99999999999; B = 27 },
    @{ Row = 29; A = 45690.99999999999; B = 27 },
    @{ Row = 30; A = 45697.99999999999; B = 28 },
    @{ Row = 31; A = 45704.99999999999; B = 28 },
    @{ Row = 32; A = 45711.99999999999; B = 29 },
    @{ Row = 33; A = 45718.99999999999; B = 29 },
    @{ Row = 34; A = 45725.99999999999; B = 29 }
)

foreach ($r in $forecastTail) {
    $wsForecast.Cells.Item($r.Row, 1).NumberFormat = $wsForecast.Cells.Item(2, 1).NumberFormat
    $wsForecast.Cells.Item($r.Row, 1).Value = $r.A
    $wsForecast.Cells.Item($r.Row, 2).Value = $r.B
}
